# Updates the cryptos list (prices / volume / a few rank-order swaps)
# Mirrors the nightly "Updated cryptos list ... with GitHub Actions" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$text) {
    # Force text storage (no numeric/date coercion) while keeping the cell's
    # style index untouched (so no stray numFmt/quotePrefix style is created).
    $ws.Range($addr).Formula = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

function Set-PctCell([string]$addr, [string]$text) {
    # Percent-change cells already start/end with spaces, so Excel never
    # treats them as numeric - a plain value assignment is enough.
    $ws.Range($addr).Value = $text
}

# Row 2 - Bitcoin
Set-TextCell "D2" "30.655.63"
Set-PctCell  "E2" "  +0.50%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.968.22"
Set-PctCell  "E3" "  +2.79%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "0.9991"
Set-PctCell  "E4" "  -0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "248.77"
Set-PctCell  "E5" "  +1.44%  "

# Row 6 - USDC
Set-TextCell "D6" "0.9993"
Set-PctCell  "E6" "  +0.00%  "

# Row 7 - XRP
Set-TextCell "D7" "0.4816"
Set-PctCell  "E7" "  -0.02%  "

# Row 8 - Cardano
Set-TextCell "D8" "0.2948"
Set-PctCell  "E8" "  +1.96%  "

# Row 9 - Dogecoin
Set-TextCell "D9" "0.06818"
Set-PctCell  "E9" "  +1.50%  "

# Row 10 - Litecoin
Set-TextCell "D10" "111.80"
Set-PctCell  "E10" "  +0.46%  "

# Row 11 - Solana (price unchanged, only volume)
Set-PctCell  "E11" "  +1.11%  "

# Row 12 - WrappedEther
Set-TextCell "D12" "1.951.59"
Set-PctCell  "E12" "  +1.88%  "

# Row 13 - TRON
Set-TextCell "D13" "0.07717"
Set-PctCell  "E13" "  +2.16%  "

# Row 14 - Polkadot
Set-TextCell "D14" "5.492"
Set-PctCell  "E14" "  +4.38%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.6902"
Set-PctCell  "E15" "  +3.04%  "

# Row 16 - BitcoinCash
Set-TextCell "D16" "295.56"
Set-PctCell  "E16" "  +2.90%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "30.655.41"
Set-PctCell  "E17" "  +0.55%  "

# Row 18 - Avalanche (price unchanged, only volume)
Set-PctCell  "E18" "  +3.10%  "

# Row 19 - Uniswap
Set-TextCell "D19" "5.680"
Set-PctCell  "E19" "  +3.85%  "

# Row 20 / 21 swap ranks: WrappedliquidstakedEther2.0 <-> ShibaInu
Set-TextCell "B20" "ShibaInu"
Set-TextCell "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D20" "0.000007680"
Set-PctCell  "E20" "  +0.70%  "

Set-TextCell "B21" "WrappedliquidstakedEther2.0"
Set-TextCell "C21" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D21" "2.200.42"
Set-PctCell  "E21" "  +1.82%  "

# Row 22 - Dai
Set-TextCell "D22" "0.9989"
Set-PctCell  "E22" "  -0.06%  "

# Row 23 - BinanceUSD
Set-TextCell "D23" "0.9992"
Set-PctCell  "E23" "  +0.00%  "

# Row 24 - Chainlink
Set-TextCell "D24" "6.643"
Set-PctCell  "E24" "  +3.67%  "

# Row 25 - Cosmos
Set-TextCell "D25" "9.794"
Set-PctCell  "E25" "  +3.61%  "

# Row 26 - Monero
Set-TextCell "D26" "168.98"
Set-PctCell  "E26" "  +3.06%  "

# Row 27 - EthereumClassic (price unchanged, only volume)
Set-PctCell  "E27" "  -0.17%  "

# Row 28 - LidoDAOToken
Set-TextCell "D28" "2.211"
Set-PctCell  "E28" "  +3.40%  "

# Row 29 - Stellar
Set-TextCell "D29" "0.1090"
Set-PctCell  "E29" "  +2.83%  "

# Row 30 - Toncoin (price unchanged, only volume)
Set-PctCell  "E30" "  +1.82%  "

# Row 31 - Filecoin
Set-TextCell "D31" "4.685"
Set-PctCell  "E31" "  +15.65%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextCell "D32" "4.442"
Set-PctCell  "E32" "  +6.75%  "

# Row 33 - Hedera
Set-TextCell "D33" "0.05091"
Set-PctCell  "E33" "  +1.93%  "

# Row 34 - ImmutableX
Set-TextCell "D34" "0.7784"
Set-PctCell  "E34" "  +6.84%  "

# Row 35 - ARBITRUM
Set-TextCell "D35" "1.175"
Set-PctCell  "E35" "  +3.61%  "

# Row 36 - VeChain
Set-TextCell "D36" "0.02067"
Set-PctCell  "E36" "  +0.65%  "

# Row 37 - HuobiToken
Set-TextCell "D37" "2.729"
Set-PctCell  "E37" "  +0.39%  "

# Row 38 - MXToken
Set-TextCell "D38" "2.716"
Set-PctCell  "E38" "  +1.88%  "

# Row 39 - RenderToken
Set-TextCell "D39" "2.070"
Set-PctCell  "E39" "  +2.84%  "

# Row 40 - Quant
Set-TextCell "D40" "111.44"
Set-PctCell  "E40" "  +0.71%  "

# Row 41 - TheSandbox
Set-TextCell "D41" "0.4466"
Set-PctCell  "E41" "  +0.80%  "

# Row 42 - FraxShare
Set-TextCell "D42" "6.054"
Set-PctCell  "E42" "  +2.74%  "

# Row 43 - TrustWalletToken
Set-TextCell "D43" "0.8740"
Set-PctCell  "E43" "  +1.10%  "

# Row 44 / 45 swap ranks: PaxDollar <-> Aave
Set-TextCell "B44" "Aave"
Set-TextCell "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D44" "70.22"
Set-PctCell  "E44" "  +3.32%  "

Set-TextCell "B45" "PaxDollar"
Set-TextCell "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D45" "1.000"
Set-PctCell  "E45" "  +0.12%  "

# Row 46 - Aptos (price unchanged, only volume)
Set-PctCell  "E46" "  +0.70%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "9.378"
Set-PctCell  "E47" "  +0.99%  "

# Row 48 - Algorand
Set-TextCell "D48" "0.1255"
Set-PctCell  "E48" "  +1.19%  "

# Row 49 - BitcoinSV
Set-TextCell "D49" "48.08"
Set-PctCell  "E49" "  -1.43%  "

# Row 50 - Elrond
Set-TextCell "D50" "35.77"
Set-PctCell  "E50" "  +2.58%  "

# Row 51 - NEARProtocol drops out, replaced by Maker
Set-TextCell "B51" "Maker"
Set-TextCell "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D51" "918.00"
Set-PctCell  "E51" "  +12.29%  "
